# Apply the edits described by the commit:
#  - update a few values on the "input" sheet (sheet1) to reflect the new
#    OSM-dump based test configuration (TestNL / NL_with_margin_from_EU_dump.osm.pbf)
#  - move the active selection on the "input" sheet to B7

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input")

# Order matters for how new shared-strings get appended, so write the
# "F2" value before "A2" to match the expected shared-string ordering.
$ws.Range("F2").Value = "NL_with_margin_from_EU_dump.osm.pbf"
$ws.Range("E2").Value = ""
$ws.Range("D2").Value = "Network based on OSM dump"
$ws.Range("A2").Value = "TestNL"

# Update the active selection / scroll position on the input sheet.
$ws.Activate()
$ws.Range("B7").Select()
